# Swap the presentation's two themes: the active "Integral" (Red Violet)
# theme becomes the standard "Office Theme", matching the target OOXML.
#
# VBA/COM RGB() packs a hex colour "RRGGBB" as R + G*256 + B*65536.
function ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$master = $p.Designs.Item(1).SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office Theme colour scheme, in MsoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $colorScheme.Item($i).RGB = ToRGB($officeThemeColors[$i - 1])
}
